$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.387.89'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +4.23%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.427.81'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +5.43%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '556.57'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.65%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.61'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +7.15%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("E8").Value = '  +2.31%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.426.23'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.42%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.74'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.89%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.150'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.39%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.16'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +12.29%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.860.59'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.52%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '62.301.46'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.13%  '

$ws.Range("E17").Value = '  +7.03%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.426.92'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +5.79%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.18'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +6.83%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '346.06'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +11.17%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.19'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.99%  '

$ws.Range("E22").Value = '  +4.23%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.05%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.12'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.31%  '

$ws.Range("E25").Value = '  +1.43%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"

$ws.Range("E27").Value = '  +14.40%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.13'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.14%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.35'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +14.22%  '

$ws.Range("E30").Value = '  +5.33%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0₃0783'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +8.30%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.44'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +10.81%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '172.19'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.14%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.44'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.11%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.396'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.74%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '377.24'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +18.69%  '

$ws.Range("E37").Value = '  +5.12%  '

$ws.Range("E38").Value = '  +11.77%  '

$ws.Range("E40").Value = '  -0.04%  '

$ws.Range("E41").Value = '  +12.04%  '

$ws.Range("E42").Value = '  +3.16%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '144.81'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +6.73%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.65'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +6.91%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '20.71'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +10.69%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.590'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.92%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0952'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.83%  '

$ws.Range("E49").Value = '  +5.02%  '

$ws.Range("E50").Value = '  +7.02%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0₆0217'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.14%  '
